# Apply the price / volume(1h) refresh captured in the commit diff.
# Numeric-looking strings (e.g. "1.00", "2.40") are prefixed with a literal
# leading apostrophe so Excel stores them as text (preserving trailing zeros)
# instead of silently coercing them to numbers, matching the original cells' type.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.081.88'
$ws.Range("E2").Value = '  +0.03%  '
$ws.Range("D3").Value = '2.303.35'
$ws.Range("E3").Value = '  -0.06%  '
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").Value = '''300.47'
$ws.Range("E5").Value = '  -0.55%  '
$ws.Range("E6").Value = '  -2.76%  '
$ws.Range("E7").Value = '  +3.63%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("E9").Value = '  -0.08%  '
$ws.Range("D10").Value = '''35.75'
$ws.Range("E10").Value = '  -1.80%  '
$ws.Range("E11").Value = '  -0.40%  '
$ws.Range("E12").Value = '  +0.22%  '
$ws.Range("D13").Value = '''18.01'
$ws.Range("E13").Value = '  -3.21%  '
$ws.Range("D14").Value = '''6.88'
$ws.Range("E14").Value = '  -1.30%  '
$ws.Range("D15").Value = '2.662.00'
$ws.Range("E15").Value = '  +0.02%  '
$ws.Range("D16").Value = '2.329.74'
$ws.Range("E16").Value = '  +2.24%  '
$ws.Range("D17").Value = '''0.788'
$ws.Range("E17").Value = '  -1.74%  '
$ws.Range("D18").Value = '42.985.22'
$ws.Range("E18").Value = '  +0.08%  '
$ws.Range("D19").Value = '''13.23'
$ws.Range("E19").Value = '  +7.33%  '
$ws.Range("E20").Value = '  +0.57%  '
$ws.Range("E21").Value = '  -2.10%  '
$ws.Range("D22").Value = '''68.47'
$ws.Range("E22").Value = '  +0.46%  '
$ws.Range("D23").Value = '''237.91'
$ws.Range("E23").Value = '  +0.53%  '
$ws.Range("D24").Value = '''2.19'
$ws.Range("E24").Value = '  -2.66%  '
$ws.Range("D25").Value = '''0.999'
$ws.Range("E25").Value = '  -0.35%  '
$ws.Range("E26").Value = '  -1.51%  '
$ws.Range("D27").Value = '''24.78'
$ws.Range("E27").Value = '  -1.00%  '
$ws.Range("D28").Value = '''167.89'
$ws.Range("E28").Value = '  -0.84%  '
$ws.Range("D29").Value = '''9.16'
$ws.Range("E30").Value = '  -6.99%  '
$ws.Range("D31").Value = '''32.74'
$ws.Range("E31").Value = '  -6.20%  '
$ws.Range("D32").Value = '''0.999'
$ws.Range("E32").Value = '  +0.04%  '
$ws.Range("E33").Value = '  +2.15%  '
$ws.Range("B34").Value = 'RenderToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D34").Value = '''4.83'
$ws.Range("E34").Value = '  +2.08%  '
$ws.Range("B35").Value = 'Celestia'
$ws.Range("C35").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D35").Value = '''18.09'
$ws.Range("E35").Value = '  +1.51%  '
$ws.Range("D36").Value = '''2.40'
$ws.Range("E36").Value = '  -0.24%  '
$ws.Range("E37").Value = '  -1.22%  '
$ws.Range("E38").Value = '  +0.33%  '
$ws.Range("D39").Value = '''1.78'
$ws.Range("E39").Value = '  -0.54%  '
$ws.Range("E40").Value = '  +1.39%  '
$ws.Range("E41").Value = '  -3.38%  '
$ws.Range("D42").Value = '2.008.13'
$ws.Range("E42").Value = '  +0.67%  '
$ws.Range("E43").Value = '  -0.76%  '
$ws.Range("E44").Value = '  -2.92%  '
$ws.Range("D45").Value = '''10.15'
$ws.Range("E45").Value = '  -1.05%  '
$ws.Range("D46").Value = '''17.25'
$ws.Range("E46").Value = '  -3.19%  '
$ws.Range("E47").Value = '  -3.40%  '
$ws.Range("D48").Value = '''54.29'
$ws.Range("E48").Value = '  -3.19%  '
$ws.Range("D49").Value = '2.527.57'
$ws.Range("E49").Value = '  -0.06%  '
$ws.Range("E50").Value = '  -1.46%  '
$ws.Range("B51").Value = 'BitcoinSV'
$ws.Range("C51").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D51").Value = '''72.40'
$ws.Range("E51").Value = '  +2.33%  '
